$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.636063098907471
$ws.Range("B1").Value = 1.972866654396057
$ws.Range("C1").Value = 2.120338201522827
$ws.Range("D1").Value = 2.443690061569214
$ws.Range("E1").Value = 3.240051746368408
